$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.209.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.07%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.910.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.59%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.63%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'368.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.48%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'103.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.35%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.541"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -3.59%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.29%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.592"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -4.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'36.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.57%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.27%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0836"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.59%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'18.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.31%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.362.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.95%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.22%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.899.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.57%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.940"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.30%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'50.873.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.71%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.54%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.99%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'12.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.69%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.0₃0946"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.26%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'68.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.59%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'260.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.75%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.93%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.99%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.12%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'25.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.62%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -6.74%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.71%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'9.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.28%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.62%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.06%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'34.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.40%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'50.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.69%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.18%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0419"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.20%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.65%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.77%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'17.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.98%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.80%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -3.23%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'22.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.62%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'119.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.65%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -3.78%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.018.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.31%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E48").Value = "'  -5.07%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.182.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.19%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.236"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.86%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0312"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -8.03%  "
$ws.Range("E51").Style = "Normal"
